# "Options " sheet: column C is an empty spacer column sitting between the
# "Engine_power" table (A) and the "Model" table (originally D:E). The user
# selected that empty column and deleted it, which shifts every column to
# its right (Model, model_no, Transmission, Transmission_no) one place to
# the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Options ")

# Remember each table's current location (row/col/size) before the shift so
# we can re-point the table definitions afterwards - deleting a column does
# not automatically resize/re-anchor existing ListObjects.
$origRanges = @{}
foreach ($lo in $ws.ListObjects) {
    $origRanges[$lo.Name] = @($lo.Range.Row, $lo.Range.Column, $lo.Range.Rows.Count, $lo.Range.Columns.Count)
}

# Select column C, then delete it entirely (Ctrl+- on a selected column).
$ws.Columns("C").Select()
$ws.Columns("C").Delete()

# Any table that started at or to the right of column C needs to move one
# column to the left to stay aligned with its (now shifted) data.
foreach ($loName in $origRanges.Keys) {
    $info = $origRanges[$loName]
    $r = $info[0]
    $c = $info[1]
    $nRows = $info[2]
    $nCols = $info[3]
    if ($c -ge 3) {
        $c = $c - 1
    }
    $topLeft = $ws.Cells.Item($r, $c)
    $bottomRight = $ws.Cells.Item($r + $nRows - 1, $c + $nCols - 1)
    $newRange = $ws.Range($topLeft, $bottomRight)
    $lo = $ws.ListObjects.Item($loName)
    $lo.Resize($newRange)
}
